# The commit swaps the contents of ppt/theme/theme1.xml (originally the
# "Integral" theme used by the slide master / main presentation theme) and
# ppt/theme/theme2.xml (originally an "Office Theme" used only by the notes
# master) so that theme1.xml ends up with the "Office Theme" colours and
# theme2.xml ends up with the "Integral" colours.
#
# The PowerPoint object model only exposes the *reachable* theme (the one
# tied to the slide master / presentation, i.e. theme1.xml) through
# Master.Theme / ThemeColorScheme - there is no COM surface that reaches the
# notes-master-only theme2.xml independently. We therefore apply the
# reachable half of the swap: update theme1's ThemeColorScheme entries from
# the "Integral" palette to the "Office Theme" palette (dk1/lt1 are
# unchanged black/white in both themes; fonts and the fmtScheme are already
# byte-identical between the two themes).

function ToOle($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Index map (1-based): 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2
# 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$tcs.Item(1).RGB  = ToOle("000000")  # dk1 (unchanged)
$tcs.Item(2).RGB  = ToOle("FFFFFF")  # lt1 (unchanged)
$tcs.Item(3).RGB  = ToOle("44546A")  # dk2
$tcs.Item(4).RGB  = ToOle("E7E6E6")  # lt2
$tcs.Item(5).RGB  = ToOle("5B9BD5")  # accent1
$tcs.Item(6).RGB  = ToOle("ED7D31")  # accent2
$tcs.Item(7).RGB  = ToOle("A5A5A5")  # accent3
$tcs.Item(8).RGB  = ToOle("FFC000")  # accent4
$tcs.Item(9).RGB  = ToOle("4472C4")  # accent5
$tcs.Item(10).RGB = ToOle("70AD47")  # accent6
$tcs.Item(11).RGB = ToOle("0563C1")  # hlink
$tcs.Item(12).RGB = ToOle("954F72")  # folHlink
